$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"

$ws.Range("D2").Value = '43.226.85'
$ws.Range("E2").Value = '  -1.73%  '
$ws.Range("D3").Value = '2.346.79'
$ws.Range("E3").Value = '  +3.87%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").Value = '0.648'
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '231.06'
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = '65.48'
$ws.Range("E7").Value = '  +1.98%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  +0.98%  '
$ws.Range("E10").Value = '  -5.13%  '
$ws.Range("D11").Value = '56.91'
$ws.Range("E11").Value = '  -0.34%  '
$ws.Range("D12").Value = '26.61'
$ws.Range("E12").Value = '  -1.55%  '
$ws.Range("D13").Value = '2.694.72'
$ws.Range("E13").Value = '  +3.78%  '
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").Value = '15.28'
$ws.Range("E15").Value = '  -2.56%  '
$ws.Range("D16").Value = '6.24'
$ws.Range("E16").Value = '  +2.52%  '
$ws.Range("D17").Value = '0.836'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").Value = '2.348.68'
$ws.Range("E18").Value = '  +3.87%  '
$ws.Range("D19").Value = '43.169.45'
$ws.Range("E19").Value = '  -1.64%  '
$ws.Range("D20").Value = '0.0₃0971'
$ws.Range("E20").Value = '  -3.11%  '
$ws.Range("D21").Value = '73.53'
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").Value = '6.17'
$ws.Range("E22").Value = '  +1.77%  '
$ws.Range("D23").Value = '247.12'
$ws.Range("E23").Value = '  -1.47%  '
$ws.Range("E24").Value = '  +20.24%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").Value = '2.44'
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  -1.62%  '
$ws.Range("D28").Value = '9.84'
$ws.Range("E28").Value = '  -2.40%  '
$ws.Range("D29").Value = '175.12'
$ws.Range("E29").Value = '  +2.48%  '
$ws.Range("E30").Value = '  +6.46%  '
$ws.Range("E31").Value = '  +8.82%  '
$ws.Range("E32").Value = '  -8.02%  '
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("E34").Value = '  +3.92%  '
$ws.Range("D35").Value = '0.0686'
$ws.Range("E35").Value = '  -2.44%  '
$ws.Range("D36").Value = '4.98'
$ws.Range("E36").Value = '  +1.70%  '
$ws.Range("D37").Value = '2.49'
$ws.Range("E37").Value = '  +8.21%  '
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("E39").Value = '  -5.51%  '
$ws.Range("E40").Value = '  -2.86%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("E42").Value = '  +8.33%  '
$ws.Range("D43").Value = '17.78'
$ws.Range("E43").Value = '  +2.69%  '
$ws.Range("E44").Value = '  +7.79%  '
$ws.Range("D45").Value = '98.37'
$ws.Range("E45").Value = '  +0.50%  '
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("E47").Value = '  -3.30%  '
$ws.Range("E48").Value = '  -1.13%  '
$ws.Range("D49").Value = '1.433.76'
$ws.Range("E49").Value = '  -0.32%  '
$ws.Range("D50").Value = '2.567.72'
$ws.Range("E50").Value = '  +3.88%  '
$ws.Range("E51").Value = '  -10.06%  '
